$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1) # ALC
$ws.Range("H2").Value = 1009.8
$ws.Range("I2").Value = 1009.8
$ws.Range("K2").Value = 1009.8
$ws.Range("M2").Value = -896.8
$ws.Range("H12").Value = 676.0769
$ws.Range("I12").Value = 525.7778
$ws.Range("J12").Value = 1014.25
$ws.Range("K12").Value = 525.7778
$ws.Range("L12").Value = 1014.25
$ws.Range("M12").Value = -355.7778
$ws.Range("N12").Value = -1354.25
$ws.Range("H18").Value = 1184.0714
$ws.Range("I18").Value = 967.46155
$ws.Range("K18").Value = 967.46155
$ws.Range("M18").Value = -683.46155
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = ""
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = ""
$ws.Range("H88").Value = 2421.7727
$ws.Range("J88").Value = 2446.6191
$ws.Range("L88").Value = 2446.6191
$ws.Range("N88").Value = -3258.6191
$ws.Range("H91").Value = 2421.7727
$ws.Range("J91").Value = 2446.6191
$ws.Range("L91").Value = 2446.6191
$ws.Range("N91").Value = -5254.6191
$ws.Range("H99").Value = 55556612
$ws.Range("J99").Value = 2750
$ws.Range("L99").Value = 8250
$ws.Range("N99").Value = -11246
$ws.Range("H101").Value = 923.3333
$ws.Range("I101").Value = 892.5
$ws.Range("J101").Value = 985
$ws.Range("K101").Value = 2677.5
$ws.Range("L101").Value = 2955
$ws.Range("M101").Value = -1055.5
$ws.Range("N101").Value = -6199
$ws.Range("H115").Value = 211.33333
$ws.Range("I115").Value = 211.33333
$ws.Range("K115").Value = 633.99999
$ws.Range("M115").Value = 933.00001
$ws.Range("H127").Value = 1223.4
$ws.Range("I127").Value = 1279.25
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 3837.75
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = 1122.25
$ws.Range("N127").Value = -12920
$ws.Range("H129").Value = 3013.1428
$ws.Range("I129").Value = 3019
$ws.Range("K129").Value = 9057
$ws.Range("M129").Value = -4057
$ws.Range("H131").Value = 773.75
$ws.Range("I131").Value = 773.75
$ws.Range("K131").Value = 2321.25
$ws.Range("M131").Value = 2718.75
$ws.Range("H138").Value = 6430.1514
$ws.Range("I138").Value = 4966.4287
$ws.Range("J138").Value = 6824.231
$ws.Range("K138").Value = 14899.2861
$ws.Range("L138").Value = 20472.693
$ws.Range("M138").Value = -9759.286100000001
$ws.Range("N138").Value = -30752.693
$ws.Range("H141").Value = 649.2857
$ws.Range("I141").Value = 649.2857
$ws.Range("K141").Value = 1947.8571
$ws.Range("M141").Value = 3232.1429
$ws = $wb.Worksheets.Item(2) # ARM
$ws.Range("H3").Value = 833
$ws.Range("J3").Value = 999
$ws.Range("L3").Value = 999
$ws.Range("N3").Value = -1229
$ws.Range("H63").Value = 6567.6665
$ws.Range("J63").Value = 12000
$ws.Range("L63").Value = 12000
$ws.Range("N63").Value = -13372
$ws.Range("H66").Value = 6567.6665
$ws.Range("J66").Value = 12000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864
$ws.Range("H113").Value = 33750
$ws.Range("J113").Value = 33750
$ws.Range("L113").Value = 33750
$ws.Range("N113").Value = -42428
$ws = $wb.Worksheets.Item(3) # BSM
$ws.Range("H96").Value = 5285.6
$ws.Range("I96").Value = 5285.6
$ws.Range("K96").Value = 5285.6
$ws.Range("M96").Value = -2539.6
$ws.Range("H105").Value = 2670
$ws.Range("I105").Value = 1005
$ws.Range("K105").Value = 1005
$ws.Range("M105").Value = 742
$ws.Range("H107").Value = 43998.1
$ws.Range("I107").Value = 68830.664
$ws.Range("K107").Value = 68830.664
$ws.Range("M107").Value = -66910.664
$ws = $wb.Worksheets.Item(4) # CRP
$ws.Range("H28").Value = 22880.334
$ws.Range("J28").Value = 22880.334
$ws.Range("L28").Value = 22880.334
$ws.Range("N28").Value = -23370.334
$ws.Range("H74").Value = 45749.875
$ws.Range("J74").Value = 49999.75
$ws.Range("L74").Value = 49999.75
$ws.Range("N74").Value = -51747.75
$ws.Range("H77").Value = 45749.875
$ws.Range("J77").Value = 49999.75
$ws.Range("L77").Value = 149999.25
$ws.Range("N77").Value = -158735.25
$ws.Range("H99").Value = 1003689.8
$ws.Range("I99").Value = 4612.25
$ws.Range("K99").Value = 4612.25
$ws.Range("M99").Value = -3114.25
$ws.Range("H126").Value = 1003689.8
$ws.Range("I126").Value = 4612.25
$ws.Range("K126").Value = 13836.75
$ws.Range("M126").Value = -11366.75
$ws = $wb.Worksheets.Item(5) # CUL
$ws.Range("H46").Value = 2176.2354
$ws.Range("J46").Value = 2557
$ws.Range("L46").Value = 7671
$ws.Range("N46").Value = -7853
$ws.Range("H50").Value = 332.8
$ws.Range("I50").Value = 332.8
$ws.Range("K50").Value = 998.4000000000001
$ws.Range("M50").Value = -517.4000000000001
$ws.Range("H53").Value = 332.8
$ws.Range("I53").Value = 332.8
$ws.Range("K53").Value = 998.4000000000001
$ws.Range("M53").Value = -517.4000000000001
$ws.Range("H86").Value = 400
$ws.Range("I86").Value = 400
$ws.Range("K86").Value = 1200
$ws.Range("M86").Value = -14
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").Value = ""
$ws.Range("H89").Value = 400
$ws.Range("I89").Value = 400
$ws.Range("K89").Value = 3600
$ws.Range("M89").Value = 2328
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").Value = ""
$ws.Range("H99").Value = 7946.5
$ws.Range("I99").Value = 6428.8335
$ws.Range("K99").Value = 19286.5005
$ws.Range("M99").Value = -17040.5005
$ws.Range("H115").Value = 3586
$ws.Range("I115").Value = 1388
$ws.Range("K115").Value = 4164
$ws.Range("M115").Value = -2989
$ws = $wb.Worksheets.Item(6) # GSM
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 4000
$ws.Range("K80").Value = 4000
$ws.Range("M80").Value = -3002
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 4000
$ws.Range("K83").Value = 20000
$ws.Range("M83").Value = -15008
$ws.Range("H97").Value = 4395
$ws.Range("I97").Value = 4358.3335
$ws.Range("J97").Value = 4450
$ws.Range("K97").Value = 4358.3335
$ws.Range("L97").Value = 4450
$ws.Range("M97").Value = -3862.3335
$ws.Range("N97").Value = -5442
$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 3500
$ws.Range("K122").Value = 10500
$ws.Range("M122").Value = -8050
$ws = $wb.Worksheets.Item(7) # LTW
$ws.Range("H50").Value = 37000
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 42400
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 42400
$ws.Range("M50").Value = -9363
$ws.Range("N50").Value = -43674
$ws.Range("H54").Value = 43000
$ws.Range("J54").Value = 43000
$ws.Range("L54").Value = 43000
$ws.Range("N54").Value = -44288
$ws.Range("H68").Value = 5999.75
$ws.Range("J68").Value = 10000
$ws.Range("L68").Value = 10000
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 5999.75
$ws.Range("J71").Value = 10000
$ws.Range("L71").Value = 50000
$ws.Range("N71").Value = -57488
$ws.Range("H122").Value = 4994.7617
$ws.Range("I122").Value = 4630.7896
$ws.Range("K122").Value = 13892.3688
$ws.Range("M122").Value = -11442.3688
$ws.Range("H132").Value = 4377.5
$ws.Range("I132").Value = 4500
$ws.Range("K132").Value = 13500
$ws.Range("M132").Value = -10970
$ws.Range("H136").Value = 2452
$ws.Range("I136").Value = 1904
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5712
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -3162
$ws.Range("N136").Value = -14100
$ws = $wb.Worksheets.Item(8) # WVR
$ws.Range("H62").Value = 6198.4
$ws.Range("I62").Value = 5500.5
$ws.Range("J62").Value = 8990
$ws.Range("K62").Value = 5500.5
$ws.Range("L62").Value = 8990
$ws.Range("M62").Value = -4876.5
$ws.Range("N62").Value = -10238
$ws.Range("H65").Value = 6198.4
$ws.Range("I65").Value = 5500.5
$ws.Range("J65").Value = 8990
$ws.Range("K65").Value = 27502.5
$ws.Range("L65").Value = 44950
$ws.Range("M65").Value = -24382.5
$ws.Range("N65").Value = -51190
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").Value = ""
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").Value = ""
$ws.Range("H80").Value = 31250
$ws.Range("J80").Value = 31250
$ws.Range("L80").Value = 31250
$ws.Range("N80").Value = -33246
$ws.Range("H83").Value = 31250
$ws.Range("J83").Value = 31250
$ws.Range("L83").Value = 93750
$ws.Range("N83").Value = -103734
$ws.Range("H136").Value = 2706.0625
$ws.Range("I136").Value = 2650.5
$ws.Range("J136").Value = 3095
$ws.Range("K136").Value = 7951.5
$ws.Range("L136").Value = 9285
$ws.Range("M136").Value = -5401.5
$ws.Range("N136").Value = -14385
